# Group the "Rectangle 28" / "TextBox 29" shapes on slide 1 into a new
# group shape, matching the other logo rectangle+label pairs on the
# slide (e.g. "Group 13" containing "Rectangle 26" / "TextBox 27").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$range = $s.Shapes.Range(@("Rectangle 28", "TextBox 29"))
$grp = $range.Group()
$grp.Name = "Group 14"
